$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect before writing values
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A38)
$disclosureText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-05 for illustrative purposes only and are subject to change."
$ws.Range("A38").Value = $disclosureText

# Update Weight (D) and Percent Change (E) values for rows 2-35
$ws.Range("D2").Value = 0.0361100021313497
$ws.Range("E2").Value = 0.0007770007770007137
$ws.Range("D3").Value = 0.02047626126368326
$ws.Range("E3").Value = -0.002331908278274342
$ws.Range("D4").Value = 0.01918357573142586
$ws.Range("E4").Value = 0.0008892481810833797
$ws.Range("D5").Value = 0.03786742188613532
$ws.Range("E5").Value = -0.002802101576182214
$ws.Range("D6").Value = 0.03430667010430554
$ws.Range("E6").Value = -0.0008000000000000229
$ws.Range("D7").Value = 0.01978818934343512
$ws.Range("E7").Value = 0.002513534416086705
$ws.Range("D8").Value = 0.03693412741839068
$ws.Range("E8").Value = 0.004281738385784495
$ws.Range("D9").Value = 0.02033730011964738
$ws.Range("E9").Value = 0.001354524110529276
$ws.Range("D10").Value = 0.02570577110120381
$ws.Range("E10").Value = -0.01409009724151622
$ws.Range("D11").Value = 0.02364482021315333
$ws.Range("E11").Value = 0.005124056094929941
$ws.Range("D12").Value = 0.05694789905952281
$ws.Range("E12").Value = 0.003348481224587418
$ws.Range("D13").Value = 0.02508646556207181
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0.02778427067998422
$ws.Range("E14").Value = -0.01460403492923812
$ws.Range("D15").Value = 0.03348126947636849
$ws.Range("E15").Value = 0.001535836177474437
$ws.Range("D16").Value = 0.01976854909363563
$ws.Range("E16").Value = 0.0005006257822277593
$ws.Range("D17").Value = 0.03043218446204233
$ws.Range("E17").Value = -0.01081216997736989
$ws.Range("D18").Value = 0.04217246261491203
$ws.Range("E18").Value = -0.001608825557343097
$ws.Range("D19").Value = 0.1264004645913843
$ws.Range("E19").Value = -0.000662690523525411
$ws.Range("D20").Value = 0.008916010231702912
$ws.Range("E20").Value = 0.02250296091590998
$ws.Range("D21").Value = 0.01557221842283999
$ws.Range("E21").Value = -0.01692350027517897
$ws.Range("D22").Value = 0.01651122641779902
$ws.Range("E22").Value = -0.0127632753821102
$ws.Range("D23").Value = 0.01620325709821584
$ws.Range("E23").Value = -0.03637611530542217
$ws.Range("D24").Value = 0.02152566276660619
$ws.Range("E24").Value = -0.01447060022419244
$ws.Range("D25").Value = 0.01228398148368262
$ws.Range("E25").Value = 0.007234279354479733
$ws.Range("D26").Value = 0.04180082827779699
$ws.Range("E26").Value = -0.0008237684661429601
$ws.Range("D27").Value = 0.0239191715467161
$ws.Range("E27").Value = [double]"9.81065437062334E-05"
$ws.Range("D28").Value = 0.0457454576693441
$ws.Range("E28").Value = -0.003800475059382635
$ws.Range("D29").Value = 0.05518267421845324
$ws.Range("E29").Value = 0.003623844899438344
$ws.Range("D30").Value = 0.01317661808366053
$ws.Range("E30").Value = 0.001928020565552746
$ws.Range("D31").Value = 0.02068500906155213
$ws.Range("E31").Value = -0.002293577981651418
$ws.Range("D32").Value = 0.01375455156412399
$ws.Range("E32").Value = -0.02225312934631429
$ws.Range("D33").Value = 0.04186280984534603
$ws.Range("E33").Value = 0.002061855670103308
$ws.Range("D34").Value = 0.01643281845950859
$ws.Range("E34").Value = 0.01445086705202314
$ws.Range("E35").Value = -0.00187218002997791

# Re-protect the sheet to restore original protection state
$ws.Protect()

Write-Host "Update complete"
